$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# Add alternative-name columns for the three oil benchmark rows.
$ws1.Range("B2").Value = "brent"
$ws1.Range("B3").Value = "wti"
$ws1.Range("B4").Value = "urals"

# Normalize trailing whitespace on these two labels to a non-breaking space,
# matching the style used elsewhere in the sheet (e.g. "нефтяной кокс ").
$nbsp = [char]0x00A0
$ws1.Range("A8").Value = "дизель" + $nbsp
$ws1.Range("A45").Value = "npk-удобрения" + $nbsp
